$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the existing "sum" column (G1) onto the new H1 cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new "Save" column header and its values
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
